# Update the ObjTables header markers embedded in row 1 (and row 2 of the
# table-of-contents sheet) of every worksheet to match the new obj_tables /
# wc_lang header convention: lower-camel-case attribute names, plus new
# 'schema' and 'tableFormat' attributes.

$wb = $excel.ActiveWorkbook

# sheet name -> list of (cell, new value)
$updates = @{
    "!!_Table of contents" = @(
        @{ Cell = "A1"; Value = "!!!ObjTables objTablesVersion='0.0.8'" },
        @{ Cell = "A2"; Value = "!!ObjTables type='Schema' objTablesVersion='0.0.8' tableFormat='row'" }
    )
    "!!Model" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Model' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'" }
    )
    "!!Taxon" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Taxon' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'" }
    )
    "!!Environment" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Environment' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'" }
    )
    "!!Submodels" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Submodel' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Compartments" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Compartment' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Species types" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='SpeciesType' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Species" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Species' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Initial species concentrations" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='DistributionInitConcentration' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Observables" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Observable' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Functions" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Function' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Reactions" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Reaction' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Rate laws" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='RateLaw' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!dFBA objectives" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='DfbaObjective' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!dFBA objective reactions" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='DfbaObjReaction' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!dFBA objective species" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='DfbaObjSpecies' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Parameters" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Parameter' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Stop conditions" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='StopCondition' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Observations" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Observation' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Observation sets" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='ObservationSet' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Conclusions" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Conclusion' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!References" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Reference' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Authors" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Author' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
    "!!Changes" = @(
        @{ Cell = "A1"; Value = "!!ObjTables type='Data' id='Change' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'" }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
